$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of row 16 onto the new row 17 (border/fill/number-format
# pattern used throughout the table) before filling in the new data.
$ws.Range("A16:H16").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)
$ws.Rows.Item(17).RowHeight = 102

# New row content ("sous menu dynamique selon la page courante" entry)
$ws.Range("B17").Value = "sous menu dynamique selon la page courante"
$ws.Range("D17").Value = "controlleur imbriqué , etapp.request.get('_route')"
$ws.Range("E17").Value = [DateTime]"2015-03-06"

# H17 carries the hyperlink; set its text/link before G17 so the shared-string
# table keeps the same ordering as the source workbook.
$ws.Range("H17").Value = "http://www.developpez.net/forums"
$ws.Hyperlinks.Add($ws.Range("H17"), "http://www.developpez.net/forums")

$ws.Range("G17").Value = "reponse a ma question dans le forum"

# Restore H17's "hyperlink" cell style (Hyperlinks.Add() re-stamps the cell with a
# slightly different style record); copy it back from H16 which already uses it.
$ws.Range("H16").Copy()
$ws.Range("H17").PasteSpecial(-4122)

# Update the sheet's view to match: selection / top-left anchored near the new row.
$ws.Range("E17").Select() | Out-Null
